$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries (rows swap their displayed name) ---
# Nigeria / Japon swap (row 51 <-> row 52)
$ws.Range("A51").Value = "Japon"
$ws.Range("A52").Value = "Nigeria"

# Santa Lucia / Timor Oriental swap (row 202 <-> row 203)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 01:01"

# --- Update statistics (Estados Unidos) ---
$ws.Range("B4").Value = 5195332
$ws.Range("C4").Value = 43737
$ws.Range("D4").Value = 2661349
$ws.Range("E4").Value = 2368421
$ws.Range("G4").Value = 479
$ws.Range("H4").Value = 165562

# --- Update statistics (Brasil) ---
$ws.Range("D5").Value = 2118460
$ws.Range("E5").Value = 815913

# --- Update statistics (Colombia) ---
$ws.Range("B11").Value = 387481
$ws.Range("C11").Value = 10611
$ws.Range("D11").Value = 212688
$ws.Range("E11").Value = 161951
$ws.Range("G11").Value = 302
$ws.Range("H11").Value = 12842

# --- Update statistics (Alemania) ---
$ws.Range("B22").Value = 217281
$ws.Range("C22").Value = 385
$ws.Range("E22").Value = 10620

# --- Update statistics (Guatemala) ---
$ws.Range("B46").Value = 56605
$ws.Range("C46").Value = 416
$ws.Range("D46").Value = 44598
$ws.Range("E46").Value = 9796
$ws.Range("G46").Value = 14
$ws.Range("H46").Value = 2211

# --- Update statistics (row 51, now Japon) ---
$ws.Range("B51").Value = 46783
$ws.Range("C51").Value = 1344
$ws.Range("D51").Value = 32312
$ws.Range("E51").Value = 13431
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 1040

# --- Update statistics (row 52, now Nigeria) ---
$ws.Range("B52").Value = 46577
$ws.Range("C52").Value = 437
$ws.Range("D52").Value = 33186
$ws.Range("E52").Value = 12446
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 945

# --- Update statistics (Venezuela) ---
$ws.Range("B67").Value = 25805
$ws.Range("C67").Value = 844
$ws.Range("E67").Value = 12226
$ws.Range("H67").Value = 223

# --- Update statistics (Chequia) ---
$ws.Range("B74").Value = 18353
$ws.Range("C74").Value = 118
$ws.Range("D74").Value = 12785
$ws.Range("E74").Value = 5178

# --- Update statistics (Bulgaria) ---
$ws.Range("B81").Value = 13396
$ws.Range("C81").Value = 53
$ws.Range("D81").Value = 7772
$ws.Range("E81").Value = 5177
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 447

# --- Update statistics (Noruega) ---
$ws.Range("B86").Value = 9638
$ws.Range("C86").Value = 39
$ws.Range("E86").Value = 525

# --- Update statistics (Tayikistan) ---
$ws.Range("B93").Value = 7745
$ws.Range("C93").Value = 39
$ws.Range("E93").Value = 1199

# --- Update statistics (Uruguay) ---
$ws.Range("B140").Value = 1353
$ws.Range("C140").Value = 18
$ws.Range("D140").Value = 1125
$ws.Range("E140").Value = 191
